$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$range = $ws.Range("C2:C90")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
